$wb = $excel.ActiveWorkbook

# Sheet "展览" (worksheet 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 578
$ws1.Range("F4").Value = 1250
$ws1.Range("F5").Value = 1079
$ws1.Range("F6").Value = 14051
$ws1.Range("F7").Value = 15458
$ws1.Range("F9").Value = 42
$ws1.Range("F18").Value = 76
$ws1.Range("F20").Value = 1201
$ws1.Range("F21").Value = 127
$ws1.Range("F22").Value = 62
$ws1.Range("F23").Value = 6022
$ws1.Range("F25").Value = 1080
$ws1.Range("F26").Value = 5524
$ws1.Range("F27").Value = 70
$ws1.Range("F28").Value = 137
$ws1.Range("F29").Value = 101
$ws1.Range("F30").Value = 459

# Sheet "全部类型" (worksheet 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 578
$ws4.Range("F5").Value = 1250
$ws4.Range("F6").Value = 1079
$ws4.Range("F7").Value = 14051
$ws4.Range("F8").Value = 15458
$ws4.Range("F10").Value = 42
$ws4.Range("F19").Value = 76
$ws4.Range("F21").Value = 1201
$ws4.Range("F22").Value = 127
$ws4.Range("F23").Value = 62
$ws4.Range("F25").Value = 6022
$ws4.Range("F27").Value = 1080
$ws4.Range("F28").Value = 5524
$ws4.Range("F29").Value = 70
$ws4.Range("F30").Value = 137
$ws4.Range("F31").Value = 101
$ws4.Range("F32").Value = 459
